# Apply the StructureDefinition-tenant-id.xlsx edit:
#  - Metadata sheet: update URL, Version, Date, Publisher
#  - Elements sheet: move the ele-1/ext-1 Constraint(s) text from the
#    "Extension" row (row 2) down to the "Extension.extension" row (row 4)

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/tenant-id"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" root element; its Constraint(s) column (AI) loses
# the ele-1/ext-1 invariant text.
$elements.Range("AI2").Value = ""

# Row 4 = "Extension.extension" element; it gains that same invariant
# text in its Constraint(s) column (AI).
$elements.Range("AI4").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
